# Apply updated "想去人数" (column F) values to the 展览 sheet and the
# 全部类型 sheet, matching the data refresh captured in the target diff.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (row numbers as they appear on that sheet) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$exhibitUpdates = @{
    3  = 1020
    4  = 13434
    6  = 1014
    7  = 9
    8  = 1731
    10 = 119
    11 = 76
    13 = 30
    14 = 13416
    16 = 592
    17 = 8934
    18 = 6
    19 = 7998
    21 = 8
    27 = 1016
    31 = 202
    32 = 167
    33 = 372
}
foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# --- Sheet "全部类型" (row numbers as they appear on that sheet) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    4  = 1020
    5  = 13434
    7  = 1014
    8  = 9
    9  = 1731
    11 = 119
    12 = 76
    14 = 30
    15 = 13416
    17 = 592
    18 = 8934
    19 = 6
    20 = 7998
    22 = 8
    28 = 1016
    34 = 202
    35 = 167
    36 = 372
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
